# 417: Add Jonah Smith
#
# Adds a new record (Jonah Smith, OM_Key 1003) to the WMT_Extract sheet,
# plus the corresponding rolled-up summary rows on Court_Reports and
# Inst_Reports.

$wb = $excel.ActiveWorkbook

$wsExtract = $wb.Worksheets.Item("WMT_Extract")
$wsCourt   = $wb.Worksheets.Item("Court_Reports")
$wsInst    = $wb.Worksheets.Item("Inst_Reports")

# ---------------------------------------------------------------------
# WMT_Extract!A4:AO4 - new raw extract record for Jonah Smith
# ---------------------------------------------------------------------
$wsExtract.Range("A4").Value = "Farringdon"
$wsExtract.Range("B4").Value = "London"
$wsExtract.Range("C4").Value = "ND01"
$wsExtract.Range("D4").Value = "KainosLDU"
$wsExtract.Range("E4").Value = "KNS"
$wsExtract.Range("F4").Value = "WMT Team"
$wsExtract.Range("G4").Value = "WMT"
$wsExtract.Range("H4").Value = "Smith"
$wsExtract.Range("I4").Value = "Jonah"
$wsExtract.Range("J4").Value = "C"
$wsExtract.Range("K4").Value = 1003
$wsExtract.Range("L4").Value = 0
$wsExtract.Range("M4").Value = 20
$wsExtract.Range("N4").Value = 10

for ($col = 15; $col -le 40; $col++) {
    $wsExtract.Cells.Item(4, $col).Value = 0
}

$wsExtract.Range("AO4").Value = 42795.628472222219
$wsExtract.Range("AO3").Copy()
$wsExtract.Range("AO4").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Court_Reports!A4:I4 - summary record for Jonah Smith
# ---------------------------------------------------------------------
$wsCourt.Range("A4").Value = "WMT Team"
$wsCourt.Range("B4").Value = "WMT"
$wsCourt.Range("C4").Value = "Smith"
$wsCourt.Range("D4").Value = 1003
$wsCourt.Range("E4").Value = "1003|WMT|C"
$wsCourt.Range("F4").Value = 0
$wsCourt.Range("G4").Value = 0
$wsCourt.Range("H4").Value = 0

$wsCourt.Range("I4").Value = 42795.628472222219
$wsCourt.Range("I3").Copy()
$wsCourt.Range("I4").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Inst_Reports!A4:H4 - summary record for Jonah Smith
# ---------------------------------------------------------------------
$wsInst.Range("A4").Value = "WMT Team"
$wsInst.Range("B4").Value = "WMT"
$wsInst.Range("C4").Value = "Smith"
$wsInst.Range("D4").Value = 1003
$wsInst.Range("E4").Value = "1003|WMT|C"
$wsInst.Range("F4").Value = 0
$wsInst.Range("G4").Value = 0

$wsInst.Range("H4").Value = 42795.628472222219
$wsInst.Range("H3").Copy()
$wsInst.Range("H4").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Selection / active-sheet bookkeeping, matching the recorded end state:
# Court_Reports and Inst_Reports each end with their new row selected
# (full-row selection), while WMT_Extract ends up the active sheet with
# F35 selected.
# ---------------------------------------------------------------------
$wsCourt.Rows.Item(4).Select()
$wsInst.Rows.Item(4).Select()

$wsExtract.Activate()
$wsExtract.Range("F35").Select()
